$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

$cell = $t.Cell(1, 1)
$cell.Range.Text = "39 x 73" + $vtab + "  7    3" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "9|    |"

$cell = $t.Cell(1, 2)
$cell.Range.Text = "49 x 10" + $vtab + "  1    0" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "9|    |"

$cell = $t.Cell(1, 3)
$cell.Range.Text = "22 x 35" + $vtab + "  3    5" + $vtab + "  ----" + $vtab + "2|    |" + $vtab + "2|    |"

$cell = $t.Cell(2, 1)
$cell.Range.Text = "50 x 78" + $vtab + "  7    8" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "0|    |"

$cell = $t.Cell(2, 2)
$cell.Range.Text = "10 x 85" + $vtab + "  8    5" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "0|    |"

$cell = $t.Cell(2, 3)
$cell.Range.Text = "77 x 33" + $vtab + "  3    3" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "7|    |"

$cell = $t.Cell(3, 1)
$cell.Range.Text = "64 x 78" + $vtab + "  7    8" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "4|    |"

$cell = $t.Cell(3, 2)
$cell.Range.Text = "80 x 25" + $vtab + "  2    5" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "0|    |"

$cell = $t.Cell(3, 3)
$cell.Range.Text = "88 x 98" + $vtab + "  9    8" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "8|    |"

$cell = $t.Cell(4, 1)
$cell.Range.Text = "58 x 60" + $vtab + "  6    0" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "8|    |"

$cell = $t.Cell(4, 2)
$cell.Range.Text = "54 x 27" + $vtab + "  2    7" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "4|    |"

$cell = $t.Cell(4, 3)
$cell.Range.Text = "44 x 78" + $vtab + "  7    8" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "4|    |"

$cell = $t.Cell(5, 1)
$cell.Range.Text = "10 x 82" + $vtab + "  8    2" + $vtab + "  ----" + $vtab + "1|    |" + $vtab + "0|    |"

$cell = $t.Cell(5, 2)
$cell.Range.Text = "51 x 41" + $vtab + "  4    1" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "1|    |"

$cell = $t.Cell(5, 3)
$cell.Range.Text = "79 x 53" + $vtab + "  5    3" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "9|    |"
